$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 780, shifting existing rows 780-821 down to 781-822
$ws.Rows.Item(780).Insert()

# Fill in the new row's data
$ws.Cells.Item(780, 1).Value = "2026/02/11"
$ws.Cells.Item(780, 2).Value = "水"
$ws.Cells.Item(780, 3).Value = 0
$ws.Cells.Item(780, 4).Value = 136
